$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.821.62"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +7.68%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'1.815.19"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +5.26%  "
$ws.Range("E3").ClearFormats()

$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = "'250.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +4.04%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.21%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").Value = "'0.4969"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +2.40%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").Value = "'0.2798"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +8.19%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").Value = "'0.06420"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +3.57%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'1.815.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +5.32%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").Value = "'16.80"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +5.21%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'0.07172"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +3.74%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'0.6527"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +7.37%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'83.87"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +9.13%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = "'4.731"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +5.65%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = "'28.792.90"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +8.38%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.15%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'0.000007436"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +3.91%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +0.29%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").Value = "'12.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +7.64%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").Value = "'2.053.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +5.21%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = "'4.629"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +4.53%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'8.925"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +4.27%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'5.372"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +6.06%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = "'143.74"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +4.88%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").Value = "'132.57"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +25.16%  "
$ws.Range("E26").ClearFormats()

$ws.Range("D27").Value = "'16.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +8.85%  "
$ws.Range("E27").ClearFormats()

$ws.Range("E28").Value = "'  +7.63%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").Value = "'1.398"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +1.31%  "
$ws.Range("E29").ClearFormats()

$ws.Range("D30").Value = "'4.190"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +6.47%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").Value = "'0.08395"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +5.59%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").Value = "'3.875"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +4.93%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").Value = "'0.04976"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +10.70%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = "'1.095"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +8.55%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = "'0.6830"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +9.92%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").Value = "'2.711"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +4.45%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").Value = "'2.755"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +13.06%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").Value = "'2.231"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +9.39%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").Value = "'0.9705"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +4.49%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'6.085"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +7.85%  "
$ws.Range("E40").ClearFormats()

$ws.Range("D41").Value = "'0.01599"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +6.95%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +0.27%  "
$ws.Range("E42").ClearFormats()

$ws.Range("B43").Value = "'Quant"
$ws.Range("B43").ClearFormats()
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C43").ClearFormats()
$ws.Range("D43").Value = "'101.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +1.42%  "
$ws.Range("E43").ClearFormats()

$ws.Range("B44").Value = "'TheSandbox"
$ws.Range("B44").ClearFormats()
$ws.Range("C44").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").Value = "'0.4123"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +7.48%  "
$ws.Range("E44").ClearFormats()

$ws.Range("D45").Value = "'7.262"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +6.13%  "
$ws.Range("E45").ClearFormats()

$ws.Range("D46").Value = "'0.1229"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +6.23%  "
$ws.Range("E46").ClearFormats()

$ws.Range("D47").Value = "'0.05522"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +2.40%  "
$ws.Range("E47").ClearFormats()

$ws.Range("D48").Value = "'8.232"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +4.52%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").Value = "'31.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +5.77%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").Value = "'0.3650"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +8.51%  "
$ws.Range("E50").ClearFormats()

$ws.Range("E51").Value = "'  +7.00%  "
$ws.Range("E51").ClearFormats()
